# Map Maker version 2
# Adds two new worksheets (Sheet4, Sheet5) with scroll/grid calculation
# tables, and makes Sheet5 the active/selected sheet (was Sheet3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet4 : "Scroll Max / Display Height / Tiles high / Scroll / Scroll %
#           / Y / Expected" calculation table
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Sheet4"

# Header row (bold + centered)
$hdr4 = $ws4.Range("B2:H2")
$hdr4.Font.Bold = $true
$hdr4.HorizontalAlignment = -4108

# Set in the same order the original authoring session entered them, so
# new shared-string ids come out in the same order.
$ws4.Range("D2").Value = "Tiles high"
$ws4.Range("H2").Value = "Expected"
$ws4.Range("C2").Value = "Display Height"
$ws4.Range("B2").Value = "Scroll Max"
$ws4.Range("F2").Value = "Scroll %"
$ws4.Range("E2").Value = "Scroll"
$ws4.Range("G2").Value = "Y"

# Data cells (only the ones that actually hold data - keep the sheet
# sparse like the original) are centered (no bold).
$data4Cells = @("B3","C3","D3","E3","G3","H3","E4","G4","E5","G5","H5","E6","H6","E7","G7","H7","E8","G8","H8")
foreach ($addr in $data4Cells) {
    $ws4.Range($addr).HorizontalAlignment = -4108
}

$ws4.Range("B3").Value = 459
$ws4.Range("C3").Value = 361
$ws4.Range("D3").Value = 10.5
$ws4.Range("E3").Value = 0
$ws4.Range("G3").Value = 160
$ws4.Range("H3").Value = 160

$ws4.Range("E4").Value = 166
$ws4.Range("G4").Value = 320

$ws4.Range("E5").Value = 298
$ws4.Range("G5").Value = 448
$ws4.Range("H5").Formula = "=448+32"

$ws4.Range("E6").Value = 327
$ws4.Range("H6").Value = 512

$ws4.Range("E7").Value = 360
$ws4.Range("G7").Value = 512
$ws4.Range("H7").Formula = "=512+32"

$ws4.Range("E8").Value = 459
$ws4.Range("G8").Value = 640
$ws4.Range("H8").Value = 640

# Scroll % column - percentage number format, centered
$pct4 = $ws4.Range("F3:F8")
$pct4.HorizontalAlignment = -4108
$pct4.NumberFormat = "0.0%"

$ws4.Range("F3").Formula = "=E3/`$B`$3"
$ws4.Range("F4").Formula = "=E4/`$B`$3"
$ws4.Range("F5").Formula = "=E5/`$B`$3"
$ws4.Range("F6").Formula = "=E6/`$B`$3"
$ws4.Range("F7").Formula = "=E7/`$B`$3"
$ws4.Range("F8").Formula = "=E8/`$B`$3"

$ws4.Columns("C:C").ColumnWidth = 13.85546875

[void]$ws4.Range("H4").Select()

# ---------------------------------------------------------------------
# Sheet5 : "Grid / X / Y" and "Final / X / Y" padding-snap table
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "Sheet5"

$ws5.Range("B2").Value = "Grid"
$ws5.Range("F2").Value = "Final"
$ws5.Range("C2").Value = "X"
$ws5.Range("D2").Value = "Y"
$ws5.Range("G2").Value = "X"
$ws5.Range("H2").Value = "Y"

$ws5.Range("B3").Value = 0
$ws5.Range("C3").Value = 0
$ws5.Range("D3").Value = 0
$ws5.Range("G3").Value = 0
$ws5.Range("H3").Value = 0

$ws5.Range("B4").Value = 0
$ws5.Range("C4").Value = 32
$ws5.Range("D4").Value = 0
$ws5.Range("G4").Value = 64
$ws5.Range("H4").Value = 0

$ws5.Range("B5").Value = 0
$ws5.Range("C5").Value = 32
$ws5.Range("D5").Value = 32
$ws5.Range("G5").Value = 64
$ws5.Range("H5").Value = 32

$ws5.Range("B6").Value = 0
$ws5.Range("C6").Value = 0
$ws5.Range("D6").Value = 32
$ws5.Range("G6").Value = 0
$ws5.Range("H6").Value = 32

$ws5.Range("B7").Value = 1
$ws5.Range("C7").Value = 32
$ws5.Range("D7").Value = 0

$ws5.Range("B8").Value = 1
$ws5.Range("C8").Value = 64
$ws5.Range("D8").Value = 0

$ws5.Range("B9").Value = 1
$ws5.Range("C9").Value = 64
$ws5.Range("D9").Value = 32

$ws5.Range("B10").Value = 1
$ws5.Range("C10").Value = 32
$ws5.Range("D10").Value = 32

# Rows that "snapped" to the grid are highlighted in red text
$ws5.Range("C4:D5").Font.Color = 255
$ws5.Range("C7:D7").Font.Color = 255
$ws5.Range("C10:D10").Font.Color = 255

[void]$ws5.Range("F9").Select()

# Sheet5 is now the active/selected sheet (Sheet3 previously had the
# lone tabSelected flag).
$ws5.Activate()
